$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update Steps (B2) and Expected Result (D2) ---
$ws.Range("B2").Value = "Enter first name.`nEnter last name.`nEnter valid Email.`nEnter valid password.`nClick on Create Account.`nClick on Sign out.`nRegister again with the same email.`n"
$ws.Range("D2").Value = "user registered succesfully and is redirected to Activation TrackID "

# --- Row 3: update Steps/Expected/Actual, add Screenshot hyperlink ---
$ws.Range("B3").Value = "Enter registerd Email.`nEnter wrong Password.`nClick on `"Sign in`""
$ws.Range("C3").Value = "unhandled exception should be handled and user friendly message should appear instead"
$ws.Range("D3").Value = "unhandled exception appears & user friendly error message both appears"

# --- Row 4: update Steps/Expected/Actual, add Screenshot hyperlink ---
$ws.Range("B4").Value = "Enter unregisterd Email.`nEnter  Password.`nClick on `"Sign in`""
$ws.Range("C4").Value = "unhandled exception should be handled and user friendly message should appear instead"
$ws.Range("D4").Value = "unhandled exception appears & user friendly error message both appears"

# --- Row 5: update Steps/Expected/Actual, add Screenshot hyperlink ---
$ws.Range("B5").Value = "Enter first name.`nEnter last name.`nEnter invalid Email.`nEnter valid password."
$ws.Range("C5").Value = "unhandled exception should be handled and user friendly message should appear instead"
$ws.Range("D5").Value = "unhandled exception appears & user friendly error message both appears"

# --- Row 6: "No email is sent" bug moves from row 6 (unchanged content, now row 6) ---
$ws.Range("A6").Value = "No email is sent when user forgets password"
$ws.Range("B6").Value = "Register user succesfully.`nUser Signout.`nUser clicks on forget password.`nUser enters email registered.`n"
$ws.Range("C6").Value = "Email should be sent contains login information "
$ws.Range("D6").Value = "No Email is sent to user"
$ws.Range("E6").Value = "Medium"
$ws.Range("F6").Value = "High"

# --- Row 7: "exceeds limits" bug now fully filled in, taller row, new screenshot hyperlink ---
$ws.Range("A7").Value = "User gets unhandled exception with exceeds limits "
$ws.Range("B7").Value = "After Investigation , this error appears when trying to login or register in Edge, Chrome or Firefox after login with invalid data multiple times.`nIn each request , a header is sent (X-Tractive-Client) with the same value.`nThe response for login or register requests doesn't have detailed information.`nThe error disappeared after a day.`n"
$ws.Range("C7").Value = "Instead of unhandled exception , user should be provided with a clear message describe the issue and how to be able to access system again."
$ws.Range("D7").Value = "Error message doesn’t provide any information on the problem and how to fix it"
$ws.Rows.Item(7).RowHeight = 135

# --- Screenshots column: replace embedded "place in cell" images with hyperlinks to the
#     corresponding screenshot files (display text = file name) ---
$ws.Hyperlinks.Add($ws.Range("G3"), "wrongpassword.png", "", "", "wrongpassword.png")
$ws.Hyperlinks.Add($ws.Range("G4"), "notregisteredemail.png", "", "", "notregisteredemail.png")
$ws.Hyperlinks.Add($ws.Range("G5"), "unhandle2.png", "", "", "unhandle2.png")
$ws.Hyperlinks.Add($ws.Range("G7"), "response.png", "", "", "response.png")
$ws.Range("G7").Value = "response.png`n"
$ws.Range("G7").WrapText = $true

# --- Selection matches the final saved cursor position ---
$ws.Range("C7").Select()

